$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.121.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.347.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.05%  "
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.708.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.368.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.762"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.079.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  -4.96%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0988"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.964.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.19%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.566.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.75%  "
